{"js": "// Fix a template typo: the \"Other variants\" section's Footnotes line\n// wrongly reused \"fail_comments\" (copy/pasted from the \"QC failed\n// variants\" section above it) instead of \"other_comments\".\n//\n// There are two \"{% for c in fail_comments %}\" occurrences in the\n// document: one that belongs to the \"QC failed variants\" table (must\n// stay as-is) and one that belongs to the \"Other variants\" table\n// (must become \"other_comments\"). Disambiguate by walking up from each\n// match to its enclosing table and checking that table's own text.\n\nconst body = context.document.body;\n\nconst matches = body.search(\"fail_comments\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"text\");\nawait context.sync();\n\nconst tables = [];\nfor (let i = 0; i < matches.items.length; i++) {\n  const tbl = matches.items[i].parentTableOrNullObject;\n  tbl.load(\"values\");\n  tables.push(tbl);\n}\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < matches.items.length; i++) {\n  const tbl = tables[i];\n  if (tbl.isNullObject) continue;\n  const tableText = tbl.values.map((row) => row.join(\"\\t\")).join(\"\\n\");\n  if (tableText.indexOf(\"other_variants\") !== -1) {\n    target = matches.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the 'fail_comments' run inside the 'other_variants' table\");\n}\n\ntarget.insertText(\"other_comments\", \"Replace\");\nawait context.sync();\n", "ps1": "# Fix a template typo: the \"Other variants\" table's Footnotes line\n# wrongly reused \"fail_comments\" (copy/pasted from the \"QC failed\n# variants\" table above it) instead of \"other_comments\".\n#\n# The document contains two \"{% for c in fail_comments %}\" runs: one\n# that belongs to the \"QC failed variants\" table (must stay as-is) and\n# one that belongs to the \"Other variants\" table (must become\n# \"other_comments\"). Find the table whose own text mentions\n# \"other_variants\" (its `{%tr for r in other_variants %}` loop header)\n# and fix only the \"Footnotes:\" line inside that table.\n\n$d = $word.ActiveDocument\n\n$targetTable = $null\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    if ($t.Range.Text.Contains(\"other_variants\")) {\n        $targetTable = $t\n        break\n    }\n}\n\nif ($targetTable -eq $null) {\n    throw \"Could not locate the table for the 'other_variants' section\"\n}\n\n$rng = $targetTable.Range\n$f = $rng.Find\n$f.Text = \"fail_comments\"\n$f.Replacement.Text = \"other_comments\"\n$f.Forward = $true\n$f.MatchCase = $true\n$f.MatchWholeWord = $false\n$f.MatchWildcards = $false\n\n# Replace:=wdReplaceOne (1) -- only the single occurrence inside this\n# table's range, not every occurrence in the document.\n$found = $f.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 1)\n\nif (-not $found) {\n    throw \"Could not find 'fail_comments' inside the 'other_variants' table\"\n}\n"}
